$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "along"
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 5).Value = "variant_warm_glowdonation"
}
